$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Previous account-statement rows are removed and replaced with the new
# worker data (period moves from 2506 to 2507), per commit message:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"

# Row 16: CC / 1049929523 / JAIDER PATERNINA MEJIA / 2507
$ws.Range("C16").Value = "1049929523"
$ws.Range("D16").Value = "JAIDER PATERNINA MEJIA"
$ws.Range("E16").Value = "2507"

# Row 17: CC / 1128058528 / JORGE ENRIQUE TABORDA CARRILLO / 2507
$ws.Range("C17").Value = "1128058528"
$ws.Range("D17").Value = "JORGE ENRIQUE TABORDA CARRILLO"
$ws.Range("E17").Value = "2507"

# Row 18: CC / 1050954661 / KARLA LICETH ARRIETA TAPIA / 2507
$ws.Range("C18").Value = "1050954661"
$ws.Range("D18").Value = "KARLA LICETH ARRIETA TAPIA"
$ws.Range("E18").Value = "2507"
